# save data done + era data updated
# Add a new "Save" column (H) to the sheet:
#  - H1: header "Save", styled like the other header cells (copy format from G1)
#  - H2, H3: numeric 0 values (unstyled, like the other data cells)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell, matching the look of the existing headers (row 1).
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats

# New data cells for the "Save" column.
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
